# Fruta / hortaliza, semanal
# Reorders the weekly date-blocks in rows 2-28 (Femacal de La Calera - Esparragos)
# Each block of rows sharing the same "Fecha" (and its Volumen/Precio data) moves as a
# unit to a new position in the table; row 23-24 (25-Nov block) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the whole data block (rows 2-28, columns A-R) before making any changes.
$snapshot = $ws.Range("A2:R28").Value2

# 0-based source-row offsets (relative to row 2) for each destination row 2..28, in order.
# i.e. destination row (2 + $i) should receive the data currently living in
# source row (2 + $sourceOrder[$i]).
$sourceOrder = @(17,18,23,24,11,12,2,25,26,15,16,9,10,19,20,3,4,5,6,0,1,21,22,7,8,13,14)

$rowCount = 27
$colCount = 18

$newArr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcIdx = $sourceOrder[$i]
    for ($c = 0; $c -lt $colCount; $c++) {
        $newArr[$i,$c] = $snapshot[$srcIdx + 1, $c + 1]
    }
}

$ws.Range("A2:R28").Value2 = $newArr
